$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-like updates (safe to assign directly without type coercion) ---
$ws.Range("D2").Value = "71.213.49"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.866.85"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "3.863.57"
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "4.518.91"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "3.865.46"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "71.271.70"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("E22").Value = "  +4.77%  "
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  +3.89%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "3.818.20"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("E39").Value = "  +10.70%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E41").Value = "  +8.21%  "
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -6.89%  "
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E49").Value = "  +3.79%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("E51").Value = "  -2.65%  "

# --- Numeric-looking "Price" text values that must remain stored as text ---
# Excel auto-converts plain numeric-looking strings to real numbers when
# assigned via .Value, which would change cell type/formatting vs. the source
# (inline text). Force text storage by temporarily setting the cell number
# format to Text ("@"), assigning the values, then clearing the format again
# so the cells end up with no explicit style, matching the original workbook.
$numRng = $ws.Range("D2:D51")
$numRng.NumberFormat = "@"
$ws.Range("D5").Value = "700.17"
$ws.Range("D6").Value = "174.18"
$ws.Range("D9").Value = "0.525"
$ws.Range("D11").Value = "7.14"
$ws.Range("D12").Value = "0.461"
$ws.Range("D14").Value = "36.58"
$ws.Range("D19").Value = "7.24"
$ws.Range("D21").Value = "11.21"
$ws.Range("D22").Value = "499.98"
$ws.Range("D23").Value = "0.726"
$ws.Range("D24").Value = "85.03"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("D26").Value = "10.80"
$ws.Range("D27").Value = "12.29"
$ws.Range("D29").Value = "3.22"
$ws.Range("D33").Value = "29.77"
$ws.Range("D34").Value = "0.182"
$ws.Range("D35").Value = "9.27"
$ws.Range("D37").Value = "0.999"
$ws.Range("D40").Value = "3.42"
$ws.Range("D41").Value = "1.05"
$ws.Range("D42").Value = "6.04"
$ws.Range("D43").Value = "0.999"
$ws.Range("D45").Value = "0.000315"
$ws.Range("D46").Value = "163.95"
$ws.Range("D49").Value = "418.03"
$ws.Range("D50").Value = "1.39"
$ws.Range("D51").Value = "43.81"
$numRng.ClearFormats()
